$d = $word.ActiveDocument

# Locate the anchor paragraph ("LOQ4084: Fenomenos de Transporte II (Requisito fraco)")
# and the trailing copyright paragraph ("(c) 2020 . Contact: luizeleno@usp.br. ...").
# Everything strictly between/after the anchor up to (and including) the copyright
# paragraph - i.e. the blank paragraph, the "Ver no Jupiter ..." paragraph, and the
# copyright paragraph itself - is removed, while the paragraphs that follow (a blank
# paragraph and the trailing page-break paragraph) are left untouched.

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*LOQ4084*Requisito fraco*") {
        $startIndex = $i + 1
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endIndex = $i
    }
}

if ($startIndex -ge 1 -and $endIndex -ge $startIndex) {
    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
